$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineShapeInHeaderFooter($hfRange, $newName) {
    $shp = $hfRange.Range.InlineShapes.Item(1)
    $shp.Select()
    $sel = $word.Selection.InlineShapes.Item(1)
    $sel.Name = $newName
}

# Footers: Item(1) = default footer (footer2.xml, docPr id=4, PearsonLogo)
#          Item(2) = first-page footer (footer1.xml, docPr id=2, PearsonLogo)
Rename-InlineShapeInHeaderFooter $sec.Footers.Item(1) "image2.png"
Rename-InlineShapeInHeaderFooter $sec.Footers.Item(2) "image2.png"

# Headers: Item(1) = default header (header2.xml, docPr id=3, BTec_Logo-Orange)
#          Item(2) = first-page header (header1.xml, docPr id=1, BTec_Logo-Orange)
Rename-InlineShapeInHeaderFooter $sec.Headers.Item(1) "image1.jpg"
Rename-InlineShapeInHeaderFooter $sec.Headers.Item(2) "image1.jpg"

Write-Output "done"
